# Row-by-row refresh of the screened-stock metrics (cols C,D,I,J,K,L,M).
# Source data pulled fresh; a couple of rows (D6/D31) move off the old
# ad-hoc "0%" style onto the sheet-standard 0.00% percent format used
# everywhere else in column D.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 24600.0
$ws.Range("D2").Value = 0.0314
$ws.Range("I2").Value = 4.07
$ws.Range("J2").Value = 87.0
$ws.Range("K2").Value = 87.0
$ws.Range("L2").Value = 100.0
$ws.Range("M2").Value = 100.0

# Row 3
$ws.Range("C3").Value = 110100.0
$ws.Range("D3").Value = -0.017
$ws.Range("I3").Value = 5.9
$ws.Range("J3").Value = 68.0
$ws.Range("K3").Value = 68.0
$ws.Range("L3").Value = 38.0
$ws.Range("M3").Value = 38.0

# Row 4
$ws.Range("C4").Value = 488000.0
$ws.Range("D4").Value = -0.0111
$ws.Range("I4").Value = 3.89
$ws.Range("J4").Value = 89.0
$ws.Range("K4").Value = 89.0
$ws.Range("L4").Value = 81.0
$ws.Range("M4").Value = 81.0

# Row 5
$ws.Range("C5").Value = 29550.0
$ws.Range("D5").Value = -0.0166
$ws.Range("I5").Value = 6.77
$ws.Range("J5").Value = 35.0
$ws.Range("K5").Value = 35.0
$ws.Range("L5").Value = 20.0
$ws.Range("M5").Value = 20.0

# Row 6
$ws.Range("D6").NumberFormat = "0.00%"
$ws.Range("C6").Value = 26200.0
$ws.Range("D6").Value = -0.0038
$ws.Range("I6").Value = 4.58
$ws.Range("J6").Value = 61.0
$ws.Range("K6").Value = 61.0
$ws.Range("L6").Value = 40.0
$ws.Range("M6").Value = 40.0

# Row 7
$ws.Range("C7").Value = 27800.0
$ws.Range("D7").Value = -0.0263
$ws.Range("I7").Value = 4.32
$ws.Range("J7").Value = 79.0
$ws.Range("K7").Value = 79.0
$ws.Range("L7").Value = 44.0
$ws.Range("M7").Value = 44.0

# Row 8
$ws.Range("C8").Value = 82100.0
$ws.Range("D8").Value = -0.0509
$ws.Range("I8").Value = 3.78
$ws.Range("J8").Value = 37.0
$ws.Range("K8").Value = 37.0
$ws.Range("L8").Value = 63.0
$ws.Range("M8").Value = 63.0

# Row 9
$ws.Range("C9").Value = 11200.0
$ws.Range("D9").Value = 0.0063
$ws.Range("I9").Value = 4.6
$ws.Range("J9").Value = 93.0
$ws.Range("K9").Value = 93.0
$ws.Range("L9").Value = 98.0
$ws.Range("M9").Value = 98.0

# Row 10
$ws.Range("C10").Value = 133900.0
$ws.Range("D10").Value = 0.0167
$ws.Range("I10").Value = 2.24
$ws.Range("J10").Value = 76.0
$ws.Range("K10").Value = 76.0
$ws.Range("L10").Value = 60.0
$ws.Range("M10").Value = 60.0

# Row 11
$ws.Range("C11").Value = 264000.0
$ws.Range("D11").Value = -0.0186
$ws.Range("I11").Value = 4.55
$ws.Range("J11").Value = 76.0
$ws.Range("K11").Value = 76.0
$ws.Range("L11").Value = 60.0
$ws.Range("M11").Value = 60.0

# Row 12
$ws.Range("C12").Value = 138800.0
$ws.Range("D12").Value = -0.0014
$ws.Range("I12").Value = 4.9
$ws.Range("J12").Value = 91.0
$ws.Range("K12").Value = 91.0
$ws.Range("L12").Value = 87.0
$ws.Range("M12").Value = 87.0

# Row 13
$ws.Range("C13").Value = 19930.0
$ws.Range("D13").Value = -0.0325
$ws.Range("I13").Value = 4.77
$ws.Range("J13").Value = 77.0
$ws.Range("K13").Value = 77.0
$ws.Range("L13").Value = 38.0
$ws.Range("M13").Value = 38.0

# Row 14
$ws.Range("D14").Value = -0.0292
$ws.Range("J14").Value = 69.0
$ws.Range("K14").Value = 69.0
$ws.Range("L14").Value = 51.0
$ws.Range("M14").Value = 51.0

# Row 15
$ws.Range("C15").Value = 146800.0
$ws.Range("D15").Value = -0.0644
$ws.Range("I15").Value = 3.64
$ws.Range("J15").Value = 80.0
$ws.Range("K15").Value = 80.0
$ws.Range("L15").Value = 38.0
$ws.Range("M15").Value = 38.0

# Row 16
$ws.Range("C16").Value = 73500.0
$ws.Range("D16").Value = -0.0303
$ws.Range("I16").Value = 4.76
$ws.Range("J16").Value = 85.0
$ws.Range("K16").Value = 85.0
$ws.Range("L16").Value = 28.0
$ws.Range("M16").Value = 28.0

# Row 17
$ws.Range("C17").Value = 53400.0
$ws.Range("D17").Value = -0.0019
$ws.Range("I17").Value = 6.63
$ws.Range("J17").Value = 57.0
$ws.Range("K17").Value = 57.0
$ws.Range("L17").Value = 59.0
$ws.Range("M17").Value = 59.0

# Row 18
$ws.Range("C18").Value = 71200.0
$ws.Range("D18").Value = -0.0166
$ws.Range("I18").Value = 7.72
$ws.Range("J18").Value = 71.0
$ws.Range("K18").Value = 71.0
$ws.Range("L18").Value = 4.0
$ws.Range("M18").Value = 4.0

# Row 19
$ws.Range("C19").Value = 38050.0
$ws.Range("D19").Value = -0.0404
$ws.Range("I19").Value = 3.68
$ws.Range("J19").Value = 84.0
$ws.Range("K19").Value = 84.0
$ws.Range("L19").Value = 59.0
$ws.Range("M19").Value = 59.0

# Row 20
$ws.Range("C20").Value = 19700.0
$ws.Range("D20").Value = -0.0199
$ws.Range("I20").Value = 5.41
$ws.Range("J20").Value = 81.0
$ws.Range("K20").Value = 81.0
$ws.Range("L20").Value = 62.0
$ws.Range("M20").Value = 62.0

# Row 21
$ws.Range("C21").Value = 51100.0
$ws.Range("D21").Value = 0.002
$ws.Range("I21").Value = 5.48
$ws.Range("J21").Value = 73.0
$ws.Range("K21").Value = 73.0
$ws.Range("L21").Value = 69.0
$ws.Range("M21").Value = 69.0

# Row 22
$ws.Range("C22").Value = 21000.0
$ws.Range("D22").Value = -0.0094
$ws.Range("I22").Value = 5.86
$ws.Range("J22").Value = 45.0
$ws.Range("K22").Value = 45.0
$ws.Range("L22").Value = 69.0
$ws.Range("M22").Value = 69.0

# Row 23
$ws.Range("C23").Value = 48250.0
$ws.Range("D23").Value = -0.0223
$ws.Range("I23").Value = 4.15
$ws.Range("J23").Value = 69.0
$ws.Range("K23").Value = 69.0
$ws.Range("L23").Value = 6.0
$ws.Range("M23").Value = 6.0

# Row 24
$ws.Range("C24").Value = 15710.0
$ws.Range("D24").Value = 0.0208
$ws.Range("I24").Value = 4.14
$ws.Range("J24").Value = 93.0
$ws.Range("K24").Value = 93.0
$ws.Range("L24").Value = 85.0
$ws.Range("M24").Value = 85.0

# Row 25
$ws.Range("C25").Value = 152000.0
$ws.Range("D25").Value = -0.0256
$ws.Range("I25").Value = 2.96
$ws.Range("J25").Value = 84.0
$ws.Range("K25").Value = 84.0
$ws.Range("L25").Value = 20.0
$ws.Range("M25").Value = 20.0

# Row 26
$ws.Range("C26").Value = 135800.0
$ws.Range("D26").Value = -0.0138
$ws.Range("I26").Value = 3.98
$ws.Range("J26").Value = 81.0
$ws.Range("K26").Value = 81.0
$ws.Range("L26").Value = 33.0
$ws.Range("M26").Value = 33.0

# Row 27
$ws.Range("C27").Value = 16840.0
$ws.Range("D27").Value = -0.0065
$ws.Range("I27").Value = 6.95
$ws.Range("L27").Value = 29.0
$ws.Range("M27").Value = 29.0

# Row 28
$ws.Range("C28").Value = 41500.0
$ws.Range("D28").Value = -0.0166
$ws.Range("I28").Value = 3.51
$ws.Range("L28").Value = 65.0
$ws.Range("M28").Value = 65.0

# Row 29
$ws.Range("C29").Value = 281500.0
$ws.Range("D29").Value = -0.0326
$ws.Range("I29").Value = 2.66
$ws.Range("J29").Value = 84.0
$ws.Range("K29").Value = 84.0
$ws.Range("L29").Value = 52.0
$ws.Range("M29").Value = 52.0

# Row 30
$ws.Range("C30").Value = 49200.0
$ws.Range("D30").Value = -0.014
$ws.Range("I30").Value = 5.02
$ws.Range("J30").Value = 58.0
$ws.Range("K30").Value = 58.0
$ws.Range("L30").Value = 27.0
$ws.Range("M30").Value = 27.0

# Row 31
$ws.Range("D31").NumberFormat = "0.00%"
$ws.Range("C31").Value = 77400.0
$ws.Range("D31").Value = -0.0227
$ws.Range("I31").Value = 2.79
$ws.Range("J31").Value = 93.0
$ws.Range("K31").Value = 93.0
$ws.Range("L31").Value = 72.0
$ws.Range("M31").Value = 72.0

# Row 32
$ws.Range("C32").Value = 54900.0
$ws.Range("D32").Value = -0.0469
$ws.Range("I32").Value = 1.09
$ws.Range("J32").Value = 83.0
$ws.Range("K32").Value = 83.0
$ws.Range("L32").Value = 40.0
$ws.Range("M32").Value = 40.0

# Row 33
$ws.Range("C33").Value = 163000.0
$ws.Range("D33").Value = -0.0468
$ws.Range("I33").Value = 2.44
$ws.Range("J33").Value = 83.0
$ws.Range("K33").Value = 83.0
$ws.Range("L33").Value = 51.0
$ws.Range("M33").Value = 51.0

# Row 34
$ws.Range("C34").Value = 50900.0
$ws.Range("D34").Value = 0.0059
$ws.Range("I34").Value = 5.3
$ws.Range("J34").Value = 71.0
$ws.Range("K34").Value = 71.0
$ws.Range("L34").Value = 92.0
$ws.Range("M34").Value = 92.0

# Row 35
$ws.Range("C35").Value = 92000.0
$ws.Range("D35").Value = -0.0108
$ws.Range("I35").Value = 3.91
$ws.Range("J35").Value = 92.0
$ws.Range("K35").Value = 92.0
$ws.Range("L35").Value = 75.0
$ws.Range("M35").Value = 75.0

# Row 36
$ws.Range("C36").Value = 123700.0
$ws.Range("D36").Value = -0.0128
$ws.Range("I36").Value = 2.57
$ws.Range("J36").Value = 94.0
$ws.Range("K36").Value = 94.0
$ws.Range("L36").Value = 74.0
$ws.Range("M36").Value = 74.0

# Row 37
$ws.Range("C37").Value = 14690.0
$ws.Range("D37").Value = -0.0174
$ws.Range("I37").Value = 4.42
$ws.Range("J37").Value = 87.0
$ws.Range("K37").Value = 87.0
$ws.Range("L37").Value = 68.0
$ws.Range("M37").Value = 68.0

# Row 38
$ws.Range("C38").Value = 13530.0
$ws.Range("D38").Value = -0.0377
$ws.Range("I38").Value = 3.7
$ws.Range("J38").Value = 78.0
$ws.Range("K38").Value = 78.0
$ws.Range("L38").Value = 44.0
$ws.Range("M38").Value = 44.0

# Row 39
$ws.Range("C39").Value = 23050.0
$ws.Range("D39").Value = -0.0295
$ws.Range("I39").Value = 4.32
$ws.Range("J39").Value = 85.0
$ws.Range("K39").Value = 85.0
$ws.Range("L39").Value = 55.0
$ws.Range("M39").Value = 55.0

# Row 40
$ws.Range("C40").Value = 25700.0
$ws.Range("D40").Value = -0.0338
$ws.Range("I40").Value = 4.67
$ws.Range("J40").Value = 92.0
$ws.Range("K40").Value = 92.0
$ws.Range("L40").Value = 46.0
$ws.Range("M40").Value = 46.0

# Leave the cursor parked on F12, matching the saved selection state.
$ws.Range("F12").Select()
